# Managing and Running Python Effectively - targeted edits
#
# 1) Footer "date updated automatically" fields bump from 11/11/2021 to
#    12/11/2021 (slide layouts 9, 10, 11 + the notes master).
# 2) Slide 24 ("__name__"): the trailing literal "__main__" in the second
#    bullet gets its own run in the Courier New code font (matching how
#    "__name__" is styled earlier in the same text box).
# 3) Slide 27 ("unittest"): a new explanatory bullet is appended -
#    "Can be used to check code and work and that it continues to work".

$p = $ppt.ActivePresentation

# --- 1) Date placeholders on the slide layouts ------------------------
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
foreach ($layoutIndex in 9, 10, 11) {
    $layout = $layouts.Item($layoutIndex)
    foreach ($shp in $layout.Shapes) {
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "12/11/2021"
        }
    }
}

# Notes master date placeholder (best effort - some hosts don't allow
# editing notes-master shapes directly).
try {
    $notesMaster = $p.NotesMaster
    foreach ($shp in $notesMaster.Shapes) {
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "12/11/2021"
        }
    }
} catch {
    Write-Output "NotesMaster date placeholder could not be updated: $_"
}

# --- 2) Slide 24: split "__main__" into its own Courier New run -------
$slide24 = $p.Slides.Item(24)
$contentShape = $slide24.Shapes.Item(2)
$tr24 = $contentShape.TextFrame.TextRange
$fullText24 = $tr24.Text
$target = "__main__"
$firstIdx = $fullText24.IndexOf("this will have the value __main__")
$mainIdx = $fullText24.IndexOf($target, $firstIdx)
$mainRun = $tr24.Characters($mainIdx + 1, $target.Length)
$mainRun.Font.NameAscii = "Courier New"
$mainRun.Font.NameComplexScript = "Courier New"

# --- 3) Slide 27: append the new bullet paragraph ----------------------
$slide27 = $p.Slides.Item(27)
$contentShape27 = $slide27.Shapes.Item(2)
$tr27 = $contentShape27.TextFrame.TextRange
$tr27.InsertAfter("`rCan be used to check code and work and that it continues to work")

$fullText27 = $tr27.Text
$secondPart = "it continues to work"
$secondIdx = $fullText27.IndexOf($secondPart, $fullText27.IndexOf("Can be used to check code"))
$secondRun = $tr27.Characters($secondIdx + 1, $secondPart.Length)
$secondRun.Font.NameAscii = $secondRun.Font.NameAscii
